$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number but must remain text
# (matches the original inlineStr cell type / "Price" column formatting).
$textCells = @("D5", "D6", "D7", "D10", "D17", "D20", "D21", "D23", "D24", "D25", "D27", "D31", "D32", "D36", "D37", "D38", "D39", "D42", "D43", "D44", "D45", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.056.04"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "3.010.40"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "593.60"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "147.01"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.009.38"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").Value = "6.34"
$ws.Range("E10").Value = "  +8.92%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("E15").Value = "  +2.40%  "
$ws.Range("D16").Value = "3.509.84"
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "6.98"
$ws.Range("E17").Value = "  -1.84%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "62.024.11"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").Value = "3.004.69"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "446.01"
$ws.Range("D21").Value = "14.18"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").Value = "7.41"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").Value = "82.15"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").Value = "10.87"
$ws.Range("E25").Value = "  +8.69%  "
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").Value = "12.08"
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "7.17"
$ws.Range("E31").Value = "  +2.04%  "
$ws.Range("D32").Value = "2.11"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("E33").Value = "  -2.41%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").Value = "0.0₃0847"
$ws.Range("E35").Value = "  +3.52%  "
$ws.Range("D36").Value = "1.03"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "5.84"
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("D38").Value = "50.19"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").Value = "9.10"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("E40").Value = "  -4.08%  "
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("D42").Value = "0.123"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "0.285"
$ws.Range("E43").Value = "  +5.56%  "
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").Value = "40.93"
$ws.Range("E44").Value = "  +9.04%  "
$ws.Range("D45").Value = "393.46"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").Value = "2.726.13"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").Value = "134.56"
$ws.Range("E48").Value = "  +3.93%  "
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("E51").Value = "  -1.70%  "
